$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The cell E8 previously held "Good Morning"; update it to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active selection on the sheet (matches the author's saved view state)
$ws.Range("E8").Select()
